# UsersContext.xlsx - reorganize Khaleel Mustafa's org/team codes from the
# retired "EF" codes to the new "ET" codes (Strategic Collaboration &
# Contracts Team moved from EF31/EF01 to ET71/ET01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: Mustafa, Khaleel Mohammad
$ws.Range("C9").Value = "ET71"                                           # ORG_CODE
$ws.Range("K9").Value = "ET01"                                           # GROUP_Code
$ws.Range("L9").Value = "ET01-New Venture Group"                         # GROUP_NAME
$ws.Range("O9").Value = "ET71"                                           # TEAM_Code
$ws.Range("P9").Value = "ET71-Strategic Collaboration & Contracts Team"  # TEAM_NAME
